$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.619.87"
$ws.Range("E2").Value = "  +2.05%  "

$ws.Range("D3").Value = "2.599.62"
$ws.Range("E3").Value = "  +0.48%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'581.45"
$ws.Range("E5").Value = "  +5.11%  "

$ws.Range("D6").Value = "'143.06"
$ws.Range("E6").Value = "  +2.11%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("E8").Value = "  +1.00%  "

$ws.Range("D9").Value = "2.605.99"
$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").Value = "'6.51"
$ws.Range("E10").Value = "  -2.96%  "

$ws.Range("E11").Value = "  +1.94%  "

$ws.Range("D12").Value = "'0.156"
$ws.Range("E12").Value = "  -2.72%  "

$ws.Range("D13").Value = "'0.371"
$ws.Range("E13").Value = "  +4.17%  "

$ws.Range("D14").Value = "3.060.99"
$ws.Range("E14").Value = "  +0.60%  "

$ws.Range("D15").Value = "'24.72"
$ws.Range("E15").Value = "  +7.87%  "

$ws.Range("D16").Value = "60.602.94"
$ws.Range("E16").Value = "  +2.05%  "

$ws.Range("E17").Value = "  +3.27%  "

$ws.Range("D18").Value = "2.607.24"
$ws.Range("E18").Value = "  +0.63%  "

$ws.Range("D19").Value = "'11.50"
$ws.Range("E19").Value = "  +10.67%  "

$ws.Range("D20").Value = "'4.67"
$ws.Range("E20").Value = "  +2.49%  "

$ws.Range("D21").Value = "'348.43"
$ws.Range("E21").Value = "  +2.42%  "

$ws.Range("D22").Value = "'6.91"
$ws.Range("E22").Value = "  +5.22%  "

$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("D24").Value = "'0.523"
$ws.Range("E24").Value = "  +8.80%  "

$ws.Range("D25").Value = "'63.43"
$ws.Range("E25").Value = "  +0.43%  "

$ws.Range("E26").Value = "  +0.20%  "

$ws.Range("E27").Value = "  +0.26%  "

$ws.Range("D28").Value = "'8.08"
$ws.Range("E28").Value = "  +7.63%  "

$ws.Range("D29").Value = "0.0₃0798"
$ws.Range("E29").Value = "  +3.87%  "

$ws.Range("E30").Value = "  +11.50%  "

$ws.Range("D31").Value = "'6.43"
$ws.Range("E31").Value = "  +5.26%  "

$ws.Range("D32").Value = "'0.998"
$ws.Range("E32").Value = "  +0.07%  "

$ws.Range("D33").Value = "'163.41"
$ws.Range("E33").Value = "  +3.72%  "

$ws.Range("D34").Value = "'19.48"
$ws.Range("E34").Value = "  +0.63%  "

$ws.Range("D35").Value = "'4.34"
$ws.Range("E35").Value = "  +6.07%  "

$ws.Range("D36").Value = "'0.991"
$ws.Range("E36").Value = "  +9.56%  "

$ws.Range("D37").Value = "'1.25"
$ws.Range("E37").Value = "  +7.39%  "

$ws.Range("E38").Value = "  +10.39%  "

$ws.Range("D39").Value = "'38.12"
$ws.Range("E39").Value = "  +1.59%  "

$ws.Range("D40").Value = "'3.93"
$ws.Range("E40").Value = "  +6.82%  "

$ws.Range("D41").Value = "'310.65"

$ws.Range("D42").Value = "'0.839"
$ws.Range("E42").Value = "  -0.30%  "

$ws.Range("D43").Value = "'134.55"
$ws.Range("E43").Value = "  -0.86%  "

$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.24%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.0992"
$ws.Range("E45").Value = "  +1.90%  "

$ws.Range("D46").Value = "'5.02"
$ws.Range("E46").Value = "  +10.81%  "

$ws.Range("D47").Value = "'19.81"
$ws.Range("E47").Value = "  +4.16%  "

$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "'0.0550"
$ws.Range("E48").Value = "  +3.32%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.604"
$ws.Range("E49").Value = "  +0.85%  "

$ws.Range("D50").Value = "'20.11"
$ws.Range("E50").Value = "  +8.44%  "

$ws.Range("E51").Value = "  +2.72%  "
